# Add three new market-watch test sheets (RatesMarket, FuturesMarket, CryptoMarket)
# after the existing FXMarket sheet, matching the "added test cases to market watch" commit.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# RatesMarket
# ---------------------------------------------------------------------------
$fxSheet = $wb.Worksheets.Item("FXMarket")
$rates = $wb.Worksheets.Add($null, $fxSheet)
$rates.Name = "RatesMarket"

$rates.Range("A1").Value = "Assertions"
$rates.Range("A2").Value = "U.S. 10 Year Treasury Note"
$rates.Range("A3").Value = "Germany 10 Year Government Bond"
$rates.Range("A4").Value = "Italy 10 Year Government Bond"
$rates.Range("A5").Value = "Spain 10 Year Government Bond"
$rates.Range("A6").Value = "U.K. 10 Year Gilt"
$rates.Range("A7").Value = "Japan 10 Year Government Bond"
$rates.Columns.Item(1).ColumnWidth = 30.42

# ---------------------------------------------------------------------------
# FuturesMarket
# ---------------------------------------------------------------------------
$futures = $wb.Worksheets.Add($null, $rates)
$futures.Name = "FuturesMarket"

$futures.Range("A1").Value = "Assertions"
$futures.Range("A2").Value = "E-Mini Dow Continuous Contract"
$futures.Range("A3").Value = "E-Mini S&P 500 Future Continuous Contract"
$futures.Range("A4").Value = "E-Mini Nasdaq 100 Index Continuous Contract"
$futures.Range("A5").Value = "Gold Continuous Contract"
$futures.Range("A6").Value = "Silver Continuous Contract"
$futures.Range("A7").Value = "Crude Oil WTI (NYM `$/bbl) Front Month"
$futures.Columns.Item(1).ColumnWidth = 38.75

# ---------------------------------------------------------------------------
# CryptoMarket
# ---------------------------------------------------------------------------
$crypto = $wb.Worksheets.Add($null, $futures)
$crypto.Name = "CryptoMarket"

$crypto.Range("A1").Value = "Assertions"
$crypto.Range("A2").Value = "Bitcoin USD"
$crypto.Range("A3").Value = "Ethereum USD"
$crypto.Range("A4").Value = "XRP USD"
$crypto.Range("A5").Value = "Bitcoin Cash USD"
$crypto.Range("A6").Value = "Litecoin USD"
$crypto.Range("A7").Value = "Monero USD"
$crypto.Columns.Item(1).ColumnWidth = 13.92

# Record the "next row" selection on the two most recently visited sheets
# (CryptoMarket, then back to FuturesMarket, which is left as the active tab).
$crypto.Activate()
$crypto.Range("A8").Select() | Out-Null

$futures.Activate()
$futures.Range("A8").Select() | Out-Null
